$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on B:E data range so numeric-looking strings
# (e.g. "409.14") are not auto-converted to numbers, matching the
# original inlineStr/text storage of these cells.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "61.946.95"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "3.426.90"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "409.14"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").Value = "128.62"
$ws.Range("E6").Value = "  -1.79%  "
$ws.Range("D7").Value = "0.631"
$ws.Range("E7").Value = "  +6.02%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.737"
$ws.Range("E9").Value = "  +6.80%  "
$ws.Range("D10").Value = "0.141"
$ws.Range("E10").Value = "  +2.13%  "
$ws.Range("D11").Value = "42.76"
$ws.Range("E11").Value = "  +2.42%  "
$ws.Range("B12").Value = "ShibaInu"
$ws.Range("C12").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D12").Value = "0.0000219"
$ws.Range("E12").Value = "  +46.56%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "9.13"
$ws.Range("E13").Value = "  +8.63%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").Value = "0.141"
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "21.39"
$ws.Range("E15").Value = "  +7.60%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.960.91"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "3.425.19"
$ws.Range("E17").Value = "  -2.13%  "
$ws.Range("D18").Value = "12.54"
$ws.Range("E18").Value = "  +8.23%  "
$ws.Range("E19").Value = "  +7.17%  "
$ws.Range("D20").Value = "61.895.75"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "452.17"
$ws.Range("E21").Value = "  +44.77%  "
$ws.Range("D22").Value = "92.44"
$ws.Range("E22").Value = "  +9.89%  "
$ws.Range("D23").Value = "3.21"
$ws.Range("E23").Value = "  +1.29%  "
$ws.Range("D24").Value = "12.97"
$ws.Range("E24").Value = "  +1.83%  "
$ws.Range("E25").Value = "  +2.36%  "
$ws.Range("D26").Value = "33.09"
$ws.Range("E26").Value = "  +11.52%  "
$ws.Range("D27").Value = "8.77"
$ws.Range("E27").Value = "  +7.41%  "
$ws.Range("D28").Value = "4.78"
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("D29").Value = "7.70"
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("D30").Value = "2.77"
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("D31").Value = "11.98"
$ws.Range("E31").Value = "  +5.84%  "
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("D33").Value = "43.20"
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "0.0502"
$ws.Range("E36").Value = "  +3.17%  "
$ws.Range("D37").Value = "54.30"
$ws.Range("E37").Value = "  +5.35%  "
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "0.135"
$ws.Range("E40").Value = "  +7.90%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "0.322"
$ws.Range("E41").Value = "  +2.61%  "
$ws.Range("E42").Value = "  -3.01%  "
$ws.Range("D43").Value = "141.88"
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("D44").Value = "4.27"
$ws.Range("E44").Value = "  +8.87%  "
$ws.Range("E45").Value = "  +1.01%  "
$ws.Range("D46").Value = "2.52"
$ws.Range("E46").Value = "  +13.55%  "
$ws.Range("D47").Value = "16.65"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").Value = "22.30"
$ws.Range("E48").Value = "  +4.68%  "
$ws.Range("D49").Value = "2.15"
$ws.Range("E49").Value = "  +10.12%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "3.770.54"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.137"
$ws.Range("E51").Value = "  +16.08%  "

# Remove the temporary text formatting so cell styles match the original
# (no explicit style index on data rows), while keeping values as text.
$ws.Range("B2:E51").ClearFormats()
